$wb = $excel.ActiveWorkbook

# Add one day of data (2023-04-19) to the violent crime tracker.
# Each sheet is updated in its "2023" column (J), representing the running
# year-to-date total; a few sheets also have a corrected prior column (I).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1991
$ws.Range("J3").Value = 2076
$ws.Range("I4").Value = 1756
$ws.Range("J4").Value = 469
$ws.Range("I6").Value = 8970
$ws.Range("J6").Value = 2599
$ws.Range("J7").Value = 7287

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 87
$ws.Range("J7").Value = 243

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 105
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 265

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 32
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 213
$ws.Range("J8").Value = 454
$ws.Range("J9").Value = 48
$ws.Range("J10").Value = 44
$ws.Range("J11").Value = 97
$ws.Range("J15").Value = 94
$ws.Range("J19").Value = 244
$ws.Range("J20").Value = 152
$ws.Range("J22").Value = 15
$ws.Range("J24").Value = 24
$ws.Range("J25").Value = 43
$ws.Range("J27").Value = 43
$ws.Range("J29").Value = 408
$ws.Range("J30").Value = 27
$ws.Range("J32").Value = 14
$ws.Range("I33").Value = 1142
$ws.Range("J33").Value = 305
$ws.Range("J36").Value = 109
$ws.Range("J37").Value = 243
$ws.Range("J42").Value = 279
$ws.Range("J43").Value = 73
$ws.Range("J47").Value = 67
$ws.Range("J49").Value = 43
$ws.Range("J51").Value = 99
$ws.Range("J53").Value = 67
$ws.Range("J54").Value = 149
$ws.Range("J57").Value = 38
$ws.Range("J60").Value = 46
$ws.Range("I63").Value = 201
$ws.Range("J63").Value = 34
$ws.Range("J64").Value = 48
$ws.Range("J65").Value = 187
$ws.Range("J67").Value = 265
$ws.Range("J72").Value = 27
$ws.Range("J76").Value = 108
$ws.Range("J78").Value = 94
$ws.Range("J79").Value = 220
$ws.Range("J80").Value = 16
$ws.Range("J83").Value = 175
$ws.Range("J85").Value = 340
$ws.Range("J86").Value = 41
$ws.Range("J88").Value = 78
$ws.Range("J89").Value = 78
$ws.Range("J90").Value = 82
$ws.Range("J93").Value = 33
$ws.Range("J94").Value = 56
$ws.Range("J96").Value = 81
$ws.Range("J97").Value = 48
$ws.Range("J101").Value = 7287

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 91
$ws.Range("I6").Value = 363
$ws.Range("J6").Value = 108
$ws.Range("I7").Value = 1142
$ws.Range("J7").Value = 305

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 121
$ws.Range("J3").Value = 144
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 408

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 82
$ws.Range("J3").Value = 138
$ws.Range("J4").Value = 21
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 340

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 58
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 62
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 12
$ws.Range("J4").Value = 4
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 27
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 144
$ws.Range("J3").Value = 151
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 454

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 26
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 68
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 32
